# Updated symbol list on Fri Dec 23 06:19:00 UTC 2022 with GitHub Actions
#
# Column D (Price) and column G (Hora) hold numeric-looking values that are
# stored as literal text (inlineStr) in the workbook, not as real numbers.
# A plain "$ws.Range(...).Value = '123'" assignment would make Excel's
# smart-entry parser coerce a numeric-looking string into a real number, so
# instead each cell's original formatting is stashed in a scratch cell, the
# cell is forced to Text ("@") before writing the literal, and then the
# original formatting is pasted back so no spurious style/format change is
# introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" (column D) values for the rows whose price actually changed.
$prices = @{
    2  = "246.33"
    3  = "22.03"
    4  = "5.418"
    5  = "0.05778"
    6  = "3.389"
    7  = "6.335"
    8  = "0.8102"
    9  = "0.9507"
    10 = "0.1424"
    11 = "0.07503"
    12 = "0.03185"
    13 = "0.03013"
    14 = "4.157"
    15 = "0.09414"
    16 = "0.001592"
    17 = "0.04816"
    19 = "0.006187"
    20 = "0.004118"
    21 = "0.0009985"
    23 = "3.778"
    24 = "2.235"
    25 = "0.3228"
    40 = "0.03885"
    41 = "0.006372"
    42 = "0.1076"
    43 = "0.003001"
    44 = "0.006484"
    45 = "0.00005594"
    48 = "0.1481"
    49 = "0.00002101"
}

foreach ($row in $prices.Keys) {
    $cellRef = "D" + $row
    $text = $prices[$row]

    # Stash the cell's current formatting in a scratch cell.
    $ws.Range($cellRef).Copy()
    $ws.Range("Z1").PasteSpecial(-4122)

    # Force Text format so the numeric-looking literal is written verbatim
    # instead of being auto-converted to a real number.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $text

    # Restore the original formatting and tidy up the scratch cell.
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range("Z1").Clear()
}

# Every data row's "Hora" (column G) advances from 5 to 6.
for ($row = 2; $row -le 51; $row++) {
    $cellRef = "G" + $row

    $ws.Range($cellRef).Copy()
    $ws.Range("Z1").PasteSpecial(-4122)

    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = "6"

    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range("Z1").Clear()
}
